# The workbook originally has three sheets:
#   1) NewContacts       (old "sheet1" - unused sample data, no longer needed)
#   2) CreateNewContact  (old "sheet2" - the real "Sample Test" contact-create sheet)
#   3) Sheet3            (old "sheet3" - empty placeholder sheet)
#
# The new class only needs a single sheet: keep the "CreateNewContact" data
# (it becomes the one-and-only sheet) and rename it to "Sheet3", dropping the
# other two sheets entirely.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the sheets that are no longer part of the new sample test class.
[void]$wb.Worksheets("NewContacts").Delete()
[void]$wb.Worksheets("Sheet3").Delete()

# What remains is the "CreateNewContact" sheet; rename it to "Sheet3" - it is
# now the only (and therefore active) sheet in the workbook.
$ws = $wb.Worksheets("CreateNewContact")
$ws.Name = "Sheet3"
[void]$ws.Activate()

# Match the updated column widths used by the new sample test layout.
$ws.Columns.Item(1).ColumnWidth = 16.09
$ws.Columns.Item(2).ColumnWidth = 17.42
$ws.Columns.Item(3).ColumnWidth = 17.25
$ws.Columns.Item(4).ColumnWidth = 15.59
$ws.Columns.Item(5).ColumnWidth = 17.25

# Update the active selection/tab state to match the new view.
[void]$ws.Range("C13").Select()

$wb.Save()
